$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: "Teen perk" rule renamed/updated -> "New Teen Rule" (all cells stored as literal text)
$ws.Range("A20").Value = "New Teen Rule"
$ws.Range("B20").Value = "'15"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'true"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.08"
$ws.Range("D20").Style = "Normal"

# Row 21: Adult non-member (values unchanged, just rewritten with proper types)
$ws.Range("A21").Value = "Adult non-member"
$ws.Range("B21").Value = 18
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0

# Row 22: Adult member (values unchanged, just rewritten with proper types)
$ws.Range("A22").Value = "Adult member"
$ws.Range("B22").Value = 18
$ws.Range("C22").Value = $true
$ws.Range("D22").Value = 0.1

# Row 23: Senior perk (values unchanged, just rewritten; C23 explicitly present but blank)
$ws.Range("A23").Value = "Senior perk"
$ws.Range("B23").Value = 60
$ws.Range("C23").Font.Bold = $false
$ws.Range("D23").Value = 0.2

# Row 24: New empty rule row (A24 is an explicit empty-text cell, B24:D24 blank placeholders)
$ws.Range("A24").Value = "'"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Font.Bold = $false
$ws.Range("C24").Font.Bold = $false
$ws.Range("D24").Font.Bold = $false
